$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the simulation output values (K:M) for rows 15-18.
# N column formulas (=SUM(K:M)) recalc automatically.
$ws.Range("K15").Value = 0.2304
$ws.Range("L15").Value = 0.6156
$ws.Range("M15").Value = 0.1539

$ws.Range("K16").Value = 0.3054
$ws.Range("L16").Value = 0.6946

$ws.Range("K17").Value = 0.7081
$ws.Range("L17").Value = 0.2919

$ws.Range("K18").Value = 0.9407
$ws.Range("L18").Value = 0.0593

# Update view state: scroll position and active selection.
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("K17").Select()
